# Insert a new data row at row 252 (shifts existing rows 252:282 down to 253:283),
# then populate the new row with this week's price observation for Haba
# (Femacal de La Calera), matching the weekly update described in the
# commit message "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing row 252 (and everything below it) down by one row.
$ws.Rows(252).Insert()

# Populate the newly inserted row 252 with the new weekly observation.
$ws.Cells.Item(252, 1).Value = 3
$ws.Cells.Item(252, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(252, 3).Value = 'Coquimbo'
$ws.Cells.Item(252, 4).Value = 45142
$ws.Cells.Item(252, 5).Value = 5
$ws.Cells.Item(252, 6).Value = 100112026
$ws.Cells.Item(252, 7).Value = 'Haba'
$ws.Cells.Item(252, 8).Value = 'Sin especificar'
$ws.Cells.Item(252, 9).Value = 'Primera'
$ws.Cells.Item(252, 10).Value = 40
$ws.Cells.Item(252, 11).Value = 15000
$ws.Cells.Item(252, 12).Value = 15000
$ws.Cells.Item(252, 13).Value = 15000
$ws.Cells.Item(252, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(252, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(252, 16).Value = 600
$ws.Cells.Item(252, 17).Value = 25
$ws.Cells.Item(252, 18).Value = 'Hortaliza'
